# ⚡️ Datos que faltaban hasta el 10
# Rename the existing sheet, add a "metadatos" sheet after it, populate it
# with a variable dictionary, and make "metadatos" the active/selected tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "datos"

# Insert the new sheet right after "datos" (so it becomes the 2nd tab).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadatos"

function Set-Cell($ws, $addr, $value, [switch]$Bold0) {
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Font.Name = "Calibri"
    $ws.Range($addr).Font.Size = 11
    $ws.Range($addr).Font.Color = 0
}

# Row 1 — headers
Set-Cell $ws2 "A1" "Variables"
Set-Cell $ws2 "B1" "Descripción"
Set-Cell $ws2 "C1" "Fuente"
Set-Cell $ws2 "D1" "Fecha_de_extracción"

# Row 2 — anno
Set-Cell $ws2 "A2" "anno"
Set-Cell $ws2 "B2" "Año"
Set-Cell $ws2 "C2" "…"
$ws2.Range("D2").Value = 45715
$ws2.Range("D2").Font.Name = "Calibri"
$ws2.Range("D2").Font.Size = 11
$ws2.Range("D2").Font.Color = 0
$ws2.Range("D2").NumberFormat = "m/d/yyyy"

# Row 3 — codmpio
Set-Cell $ws2 "A3" "codmpio"
Set-Cell $ws2 "B3" "Código del municipio"
Set-Cell $ws2 "C3" "…"
$ws2.Range("D3").Value = 45715
$ws2.Range("D3").Font.Name = "Calibri"
$ws2.Range("D3").Font.Size = 11
$ws2.Range("D3").Font.Color = 0
$ws2.Range("D3").NumberFormat = "m/d/yyyy"

# Row 4 — p51 (variable name cell keeps the default/no style)
$ws2.Range("A4").Value = "p51"
Set-Cell $ws2 "B4" "Variable Categórica P51 Encuesta Nacional de Calidad de Vida - ENCV "
Set-Cell $ws2 "C4" "Encuesta Nacional de Calidad de Vida - ENCV del Departamento Administrativo Nacional de Estadísticas - DANE"
$ws2.Range("D4").Value = 45715
$ws2.Range("D4").Font.Name = "Calibri"
$ws2.Range("D4").Font.Size = 11
$ws2.Range("D4").Font.Color = 0
$ws2.Range("D4").NumberFormat = "m/d/yyyy"

# Row 5 — porcentaje (default style cells for A5 / B5)
$ws2.Range("A5").Value = "porcentaje"
$ws2.Range("B5").Value = "Construcción Propia"
Set-Cell $ws2 "C5" "Encuesta Nacional de Calidad de Vida - ENCV del Departamento Administrativo Nacional de Estadísticas - DAN"

# Row 6 — Observaciones (default style)
$ws2.Range("A6").Value = "Observaciones"
$ws2.Range("B6").Value = "Se calcula agrupando los datos por código de departamento (coddepto) y categoría de P51. Luego, se suma la frecuencia ponderada (FEX_C.x) para cada grupo. Posteriormente, se calcula el porcentaje dividiendo la frecuencia ponderada de cada categoría por la suma total de la frecuencia ponderada del municipio y multiplicando por 100. Finalmente, se asigna el año correspondiente (anno)."

# Column sizing to match the authored widths
$ws2.Columns.Item(3).ColumnWidth = 10.42578125

# Selection / active-cell bookkeeping on the new sheet
$ws2.Range("D17").Select()

# "metadatos" is the tab that should be selected/active when the workbook opens
$ws2.Activate()

Write-Host "done"
